$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K is needed for the "Standard Deviation" header, widen it to match.
$ws.Columns.Item(11).ColumnWidth = 18.166666666666668

# Add the new header cell with the same formatting as the neighbouring header (J4).
$ws.Range("K4").Value = "Standard Deviation"
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move/update the active selection.
$ws.Range("J8").Select()
